$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(321, 'Fajl za rad sa stringovima. Operacije nad stringovima itd…', '2_convertIntToString5', 'Uradi obradu konverzije celobrojnog broja u string.')
    ,@(322, 'Fajl za rad sa stringovima. Operacije nad stringovima itd…', '2_convertIntToString6', 'Pretvaranje celih brojeva.')
    ,@(323, 'Fajl za rad sa stringovima. Operacije nad stringovima itd…', '2_convertIntToString7', 'Pretvaranje brojeva u stringove.')
    ,@(324, 'Fajl za rad sa stringovima. Operacije nad stringovima itd…', '2_convertIntToString8', 'Prebaci ceo broj u string.')
    ,@(325, 'Fajl za rad sa stringovima. Operacije nad stringovima itd…', '2_convertIntToString9', 'Konvertovanje u string.')
    ,@(326, 'Fajl sa razlicitim strukturama podataka, Stekovi, redovi itd…', '3_pq1', 'Ubaci u prioritetni red podatke izvucene iz zahteva.')
    ,@(327, 'Fajl sa razlicitim strukturama podataka, Stekovi, redovi itd…', '3_pq2', 'Uradi obradu zahteva po prioritetima.')
    ,@(328, 'Fajl sa razlicitim strukturama podataka, Stekovi, redovi itd…', '3_pq3', 'Varijable izvucene iz klijentskog zahteva ubaciti u prioritetni red.')
    ,@(329, 'Fajl sa razlicitim strukturama podataka, Stekovi, redovi itd…', '3_pq4', 'Izbrisati iz prioritetnog reda deset podataka.')
    ,@(330, 'Fajl sa razlicitim strukturama podataka, Stekovi, redovi itd…', '3_pq5', 'Obraditi podataka po prioritetima.')
    ,@(331, 'Fajl za rad sa datumima.', '4_sud1', 'Konverzija tipova u datum.')
    ,@(332, 'Fajl za rad sa datumima.', '4_sud2', 'Pretvori stringove primljene kao parametre zahteva u datume.')
    ,@(333, 'Fajl za rad sa datumima.', '4_sud3', 'Od niza stringova, formiraj jedan konkatenirani string i probaj da ga pretvoris u datum.')
    ,@(334, 'Fajl za rad sa datumima.', '4_sud4', 'Uradi operaciju nad parovima stringova i pretvori rezultat u datum.')
    ,@(335, 'Fajl za rad sa datumima.', '4_sud5', 'Podeli string da vise malih delova, i iz tih delova probaj da napravis datum.')
    ,@(336, 'Fajl za rad sa datotekama.', '6_saveListF1', 'Sadrzaj dve liste kombinovane sacuvaj u okviru tekstualnog fajla.')
    ,@(337, 'Fajl za rad sa datotekama.', '6_saveListF2', 'Spoji par listi prosledjenih preko klijenstkoj zahteva i prodaj da ih sacuvas u txt fajlu.')
    ,@(338, 'Fajl za rad sa datotekama.', '6_saveListF3', 'Indeksiranu listu sacuvaj u fajlu.')
    ,@(339, 'Fajl za rad sa datotekama.', '6_saveListF4', 'Listu interno implementiranu kao jednostruku usnimiti u fajl.')
    ,@(340, 'Fajl za rad sa datotekama.', '6_saveListF5', 'Listu interno implementiranu kao dvostruko ulanacanu usnimiti u text fajl.')
    ,@(341, 'Fajlovi sa razlicitim strukturama podataka.', '8_cm1', 'Ubaci u matricu implementiranu kao dvodimenzijalni niz novu vrednost.')
    ,@(342, 'Fajlovi sa razlicitim strukturama podataka.', '8_cm2', 'Iz matrice implementirane preko dvodimenzijalnog niza, uklonoiti element na poziciji 3 3.')
    ,@(343, 'Fajlovi sa razlicitim strukturama podataka.', '8_cm3', 'Iz matrice dohvatiti samo prvi red i ispisati ga na standardnom izlazu.')
    ,@(344, 'Fajlovi sa razlicitim strukturama podataka.', '8_cm4', 'Vratiti elemente koji se nalazi na dijagonali matrice.')
    ,@(345, 'Fajlovi sa razlicitim strukturama podataka.', '8_cm5', 'Pre brisanja matrice zabune ispisi njene elemente.')
    ,@(346, 'Fajlovi za enkripciju sifre.', '12_aes1', 'Iskoristiti aes enkripciju da bi se zastitili bitni podaci.')
    ,@(347, 'Fajlovi za enkripciju sifre.', '12_aes2', 'Aes enkripcija promenljivih iz dokument menadzera.')
    ,@(348, 'Fajlovi za enkripciju sifre.', '12_aes3', 'Enkriptovanje sifre sa aes algoritmom.')
    ,@(349, 'Fajlovi za enkripciju sifre.', '12_aes4', 'Dekripcija dva ulazna parametra preko aes algoritma.')
    ,@(350, 'Fajlovi za enkripciju sifre.', '12_aes5', 'Sifrovanje osetljivih podataka raznim algoritmima za enkripciju.')
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 3).Value = $item[1]
    $ws.Cells.Item($r, 4).Value = $item[2]
    $ws.Cells.Item($r, 5).Value = $item[3]
}

$win = $excel.ActiveWindow
$ws.Range("E351").Select()
$win.ScrollRow = 325
$win.ScrollColumn = 4

